# Scheduled runner update: refresh currentAveragePrice / Leve profit columns
# (H:N) across the eight crafting-profession sheets with latest market data.
#
# For each affected row we set H..N explicitly; cells that the source data
# no longer provides (no market listings / no longer profitable to compute)
# are cleared entirely by assigning $null, which drops the cell node on
# save - matching rows that lose a trailing M/N cell in the target sheet.

$wb = $excel.ActiveWorkbook

# NOTE: this COM host does not bind PowerShell named parameters on custom
# functions, so Set-RowValues takes its arguments positionally:
#   Set-RowValues <SheetName> <Row> <Values hashtable>
function Set-RowValues {
    param(
        [string]$SheetName,
        [int]$Row,
        [System.Collections.IDictionary]$Values
    )

    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($col in @("H", "I", "J", "K", "L", "M", "N")) {
        $addr = "$col$Row"
        if ($Values.Contains($col)) {
            $ws.Range($addr).Value = $Values[$col]
        }
    }
}

# ALC -----------------------------------------------------------------
Set-RowValues "ALC" 2 @{
    H = 556.6667; I = 556.6667; J = 0; K = 556.6667; L = 0;
    M = -443.6667; N = $null
}
Set-RowValues "ALC" 43 @{
    H = 1854.25; I = 1812.5; J = 1875.125; K = 1812.5; L = 1875.125;
    M = -1743.5; N = -2013.125
}
Set-RowValues "ALC" 51 @{
    H = 2000; I = 0; J = 2000; K = 0; L = 2000;
    M = $null; N = -2968
}
Set-RowValues "ALC" 69 @{
    H = 3452.5; J = 3501.3333; L = 10503.9999; N = -12251.9999
}
Set-RowValues "ALC" 72 @{
    H = 3452.5; J = 3501.3333; L = 31511.9997; N = -40247.9997
}
Set-RowValues "ALC" 112 @{
    H = 4760.65; J = 5116.919; L = 15350.757; N = -17566.757
}

# ARM -------------------------------------------------------------------
Set-RowValues "ARM" 32 @{
    H = 926716; I = 1019646.7; J = 28385.666; K = 1019646.7; L = 28385.666;
    M = -1019359.7; N = -28959.666
}
Set-RowValues "ARM" 98 @{
    H = 0; J = 0; L = 0; N = $null
}
Set-RowValues "ARM" 102 @{
    H = 2299.9167; I = 2011; J = 3166.6667; K = 2011; L = 3166.6667;
    M = -389; N = -6410.6667
}
Set-RowValues "ARM" 122 @{
    H = 79532.53999999999; I = 127377.875; J = 2980; K = 382133.625;
    L = 8940; M = -379683.625; N = -13840
}

# BSM ---------------------------------------------------------------------
Set-RowValues "BSM" 26 @{
    H = 45254.445; I = 26161.834; J = 83439.664; K = 26161.834; L = 83439.664;
    M = -25869.834; N = -84023.664
}
Set-RowValues "BSM" 99 @{
    H = 2171.25; I = 2031.875; J = 2450; K = 2031.875; L = 2450;
    M = -533.875; N = -5446
}
Set-RowValues "BSM" 107 @{
    H = 46601.547; I = 59778.59; J = 1799.6; K = 59778.59; L = 1799.6;
    M = -57858.59; N = -5639.6
}
Set-RowValues "BSM" 134 @{
    H = 5132.7856; I = 4536.8887; J = 6205.4; K = 13610.6661; L = 18616.2;
    M = -11075.6661; N = -23686.2
}

# CRP -----------------------------------------------------------------
Set-RowValues "CRP" 28 @{
    H = 67500; J = 67500; L = 67500; N = -67990
}
Set-RowValues "CRP" 33 @{
    H = 1908.25; I = 1747.2858; J = 3035; K = 1747.2858; L = 3035;
    M = -1368.2858; N = -3793
}
Set-RowValues "CRP" 103 @{
    H = 0; I = 0; K = 0; M = $null
}
Set-RowValues "CRP" 132 @{
    H = 1388.3914; I = 1214.7059; J = 1880.5; K = 3644.1177; L = 5641.5;
    M = -1114.1177; N = -10701.5
}

# CUL -----------------------------------------------------------------
Set-RowValues "CUL" 22 @{
    H = 964.64703; I = 833; J = 992.8570999999999; K = 2499; L = 2978.5713;
    M = -2330; N = -3316.5713
}
Set-RowValues "CUL" 27 @{
    H = 964.64703; I = 833; J = 992.8570999999999; K = 2499; L = 2978.5713;
    M = -2397; N = -3182.5713
}
Set-RowValues "CUL" 76 @{
    H = 3999.75
}
Set-RowValues "CUL" 79 @{
    H = 3999.75
}
Set-RowValues "CUL" 113 @{
    H = 539.2069; I = 552.94446; K = 1658.83338; M = 511.16662
}

# GSM -----------------------------------------------------------------
Set-RowValues "GSM" 39 @{
    H = 16907.334; J = 16907.334; L = 16907.334; N = -17971.334
}
Set-RowValues "GSM" 80 @{
    H = 1091053.4; I = 1802698; J = 201497.5; K = 1802698; L = 201497.5;
    M = -1801700; N = -203493.5
}
Set-RowValues "GSM" 83 @{
    H = 1091053.4; I = 1802698; J = 201497.5; K = 9013490; L = 1007487.5;
    M = -9008498; N = -1017471.5
}
Set-RowValues "GSM" 113 @{
    H = 74090.36; J = 2933.3333; L = 2933.3333; N = -7273.3333
}

# LTW -----------------------------------------------------------------
Set-RowValues "LTW" 22 @{
    H = 60501; I = 0; J = 60501; K = 0; L = 60501;
    M = $null; N = -61091
}
Set-RowValues "LTW" 27 @{
    H = 60501; I = 0; J = 60501; K = 0; L = 60501;
    M = $null; N = -60715
}
Set-RowValues "LTW" 106 @{
    H = 68123.336; J = 68123.336; L = 68123.336; N = -70647.336
}
Set-RowValues "LTW" 132 @{
    H = 2163.2205; I = 1797.6842; J = 2824.6667; K = 5393.0526;
    L = 8474.000100000001; M = -2863.0526; N = -13534.0001
}

# WVR -----------------------------------------------------------------
Set-RowValues "WVR" 132 @{
    H = 3464; I = 3579.5386; J = 3088.5; K = 10738.6158; L = 9265.5;
    M = -8208.6158; N = -14325.5
}

"updated 31 rows across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR"
